$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (Wins/Losses/Ties), copying the
# existing header formatting (bold, centered, bordered) from AC1
# onto AD1:AF1 so the new headers match the rest of the header row.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins/Losses/Ties) for every player row.
$ws.Range("AD2:AD45").Value = 89
$ws.Range("AE2:AE45").Value = 73
$ws.Range("AF2:AF45").Value = 0
